# "Simplify actions to fit new syntax"
#
# The "Mehr als eine Woche" row's Action cell (E4) used to hold two
# differently-formatted text runs - "LOAD(TwoWays), " followed by
# "JUMP(Safety)" (the second run carried an extra <charset> on its font,
# a leftover from how the rich text was authored). The new syntax joins
# both actions into a single plain-text value separated by a semicolon,
# so the cell no longer needs per-run formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$actionCell = $ws.Range("E4")

# Replace the two-run rich text with one plain string.
$actionCell.Value = "LOAD(TwoWays); JUMP(Safety)"

# Re-assert the (already-default) font on the whole cell so the engine
# drops the now-unused per-run font/style that only existed to carry the
# second run's distinct formatting.
$actionCell.Font.Name = "Calibri"
$actionCell.Font.Size = 11
$actionCell.Font.Color = 0

# Reflect the updated cursor position that came along with this edit.
$ws.Range("F11").Select()
